# Apply cell updates per the diff: refreshed crypto price/volume data for Mon Sep 11 10:28:28 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.892.08'
$ws.Range("E2").Value = '  -0.39%  '
$ws.Range("D3").Value = '1.599.56'
$ws.Range("E3").Value = '  -2.06%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '''209.18'
$ws.Range("E5").Value = '  -2.25%  '
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("D7").Value = '''0.475'
$ws.Range("E7").Value = '  -5.74%  '
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").Value = '''0.244'
$ws.Range("E8").Value = '  -3.24%  '
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = '''0.0609'
$ws.Range("E9").Value = '  -2.42%  '
$ws.Range("D10").Value = '''17.74'
$ws.Range("E10").Value = '  -4.28%  '
$ws.Range("E11").Value = '  -0.48%  '
$ws.Range("D12").Value = '1.824.51'
$ws.Range("E12").Value = '  -1.91%  '
$ws.Range("D13").Value = '1.604.37'
$ws.Range("E13").Value = '  -1.62%  '
$ws.Range("D14").Value = '''4.03'
$ws.Range("E14").Value = '  -3.92%  '
$ws.Range("D15").Value = '''0.507'
$ws.Range("E15").Value = '  -4.32%  '
$ws.Range("D16").Value = '25.874.08'
$ws.Range("E16").Value = '  -0.47%  '
$ws.Range("D17").Value = '''60.62'
$ws.Range("E17").Value = '  -1.89%  '
$ws.Range("D18").Value = '0.0₃0716'
$ws.Range("E18").Value = '  -3.79%  '
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("D20").Value = '''189.03'
$ws.Range("E20").Value = '  -0.60%  '
$ws.Range("D21").Value = '''4.16'
$ws.Range("E21").Value = '  -2.02%  '
$ws.Range("D22").Value = '''9.27'
$ws.Range("E22").Value = '  -3.15%  '
$ws.Range("D23").Value = '''5.91'
$ws.Range("E23").Value = '  -3.73%  '
$ws.Range("D25").Value = '''141.70'
$ws.Range("E25").Value = '  -1.08%  '
$ws.Range("D26").Value = '''0.128'
$ws.Range("E26").Value = '  -4.12%  '
$ws.Range("E27").Value = '  -3.76%  '
$ws.Range("D28").Value = '''6.49'
$ws.Range("E28").Value = '  -4.15%  '
$ws.Range("D29").Value = '''14.89'
$ws.Range("E29").Value = '  -2.02%  '
$ws.Range("E30").Value = '  -2.39%  '
$ws.Range("D31").Value = '''0.0467'
$ws.Range("E31").Value = '  -3.49%  '
$ws.Range("D32").Value = '''3.06'
$ws.Range("E32").Value = '  -2.87%  '
$ws.Range("D33").Value = '''2.99'
$ws.Range("E33").Value = '  -5.22%  '
$ws.Range("D34").Value = '''2.41'
$ws.Range("E34").Value = '  -0.94%  '
$ws.Range("D35").Value = '''1.45'
$ws.Range("E35").Value = '  -3.29%  '
$ws.Range("D36").Value = '1.106.02'
$ws.Range("E36").Value = '  -2.60%  '
$ws.Range("D37").Value = '''2.36'
$ws.Range("E37").Value = '  -3.01%  '
$ws.Range("D38").Value = '''0.794'
$ws.Range("E38").Value = '  -8.89%  '
$ws.Range("D39").Value = '''0.0150'
$ws.Range("E39").Value = '  -2.95%  '
$ws.Range("D40").Value = '''0.493'
$ws.Range("E40").Value = '  -6.09%  '
$ws.Range("D41").Value = '''95.34'
$ws.Range("E41").Value = '  -3.32%  '
$ws.Range("D42").Value = '1.737.05'
$ws.Range("E42").Value = '  -1.88%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '''5.06'
$ws.Range("E43").Value = '  -3.74%  '
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").Value = '''0.740'
$ws.Range("E44").Value = '  -5.16%  '
$ws.Range("D45").Value = '0.0₆0109'
$ws.Range("E45").Value = '  -4.91%  '
$ws.Range("D46").Value = '''52.87'
$ws.Range("E46").Value = '  -4.09%  '
$ws.Range("E47").Value = '  -1.91%  '
$ws.Range("D48").Value = '''0.0510'
$ws.Range("E48").Value = '  -3.82%  '
$ws.Range("E49").Value = '  -1.03%  '
$ws.Range("E50").Value = '  +0.13%  '
$ws.Range("E51").Value = '  -2.98%  '
